# ozon fixes 15.01.2026 part 2
# Replace the data block (rows 2-16) with the new dataset, fix up per-cell
# formatting that came along with the paste, drop the stale hyperlink,
# shrink the conditional-formatting range, and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1) Fix up number-formats that drift during the paste, by copying the
#    format (not the value) from a donor cell that already carries the
#    target style. This re-uses existing style records instead of
#    minting new ones.
# ---------------------------------------------------------------------

# Column D (amounts) picks up the "14"/"15" comma-style used by columns
# B/C/E instead of its original plain style, alternating by row parity.
# NOTE: the COM bridge only honours the first area of a multi-area
# (union / comma) destination range for PasteSpecial, so every
# destination below is issued as its own single contiguous Copy/Paste.
$ws.Range("E2").Copy()
$ws.Range("D2").PasteSpecial($xlPasteFormats)

$ws.Range("E3").Copy()
$ws.Range("D3:D5").PasteSpecial($xlPasteFormats)

$ws.Range("E2").Copy()
$ws.Range("D4").PasteSpecial($xlPasteFormats)

$ws.Range("E3").Copy()
$ws.Range("D7").PasteSpecial($xlPasteFormats)

$ws.Range("E2").Copy()
$ws.Range("D6").PasteSpecial($xlPasteFormats)
$ws.Range("E2").Copy()
$ws.Range("D8").PasteSpecial($xlPasteFormats)
$ws.Range("E2").Copy()
$ws.Range("D10").PasteSpecial($xlPasteFormats)

$ws.Range("E3").Copy()
$ws.Range("D9").PasteSpecial($xlPasteFormats)

$ws.Range("E3").Copy()
$ws.Range("D11").PasteSpecial($xlPasteFormats)
$ws.Range("E3").Copy()
$ws.Range("D13").PasteSpecial($xlPasteFormats)
$ws.Range("E3").Copy()
$ws.Range("D15").PasteSpecial($xlPasteFormats)

$ws.Range("E2").Copy()
$ws.Range("D12").PasteSpecial($xlPasteFormats)
$ws.Range("E2").Copy()
$ws.Range("D14").PasteSpecial($xlPasteFormats)
$ws.Range("E2").Copy()
$ws.Range("D16").PasteSpecial($xlPasteFormats)

# Rows 6, 8 and 10 also pick up the "odd-row" style in columns C and E.
$ws.Range("C3").Copy()
$ws.Range("C6").PasteSpecial($xlPasteFormats)
$ws.Range("C3").Copy()
$ws.Range("C8").PasteSpecial($xlPasteFormats)
$ws.Range("C3").Copy()
$ws.Range("C10").PasteSpecial($xlPasteFormats)

$ws.Range("E3").Copy()
$ws.Range("E6").PasteSpecial($xlPasteFormats)
$ws.Range("E3").Copy()
$ws.Range("E8").PasteSpecial($xlPasteFormats)
$ws.Range("E3").Copy()
$ws.Range("E10").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

# Column A (ids) for the new rows 3-10 loses its banded style entirely.
$ws.Range("A3:A10").Style = "Normal"

# ---------------------------------------------------------------------
# 2) Write the new dataset into rows 2-10.
# ---------------------------------------------------------------------

$ws.Range("A2").Value = 2264511
$ws.Range("B2").Value = 20121
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 5.5

$ws.Range("A3").Value = 2289214
$ws.Range("B3").Value = 56790
$ws.Range("C3").Value = 7
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 7.5

$ws.Range("A4").Value = 2336118
$ws.Range("B4").Value = 32413
$ws.Range("C4").Value = 12.5
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 13

$ws.Range("A5").Value = 2341441
$ws.Range("B5").Value = 15813
$ws.Range("C5").Value = 7
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 7.5

$ws.Range("A6").Value = 2347983
$ws.Range("B6").Value = 28289
$ws.Range("C6").Value = 7
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 7.5

$ws.Range("A7").Value = 2348446
$ws.Range("B7").Value = 11582
$ws.Range("C7").Value = 7
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 7.5

$ws.Range("A8").Value = 2357313
$ws.Range("B8").Value = 11036
$ws.Range("C8").Value = 7
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 7.5

$ws.Range("A9").Value = 2362348
$ws.Range("B9").Value = 107252
$ws.Range("C9").Value = 7
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 7.5

$ws.Range("A10").Value = 2367671
$ws.Range("B10").Value = 47751
$ws.Range("C10").Value = 7
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 7.5

# ---------------------------------------------------------------------
# 3) Rows 11-16 no longer carry data (only their existing banded
#    formatting remains).
# ---------------------------------------------------------------------

$ws.Range("A11:E16").ClearContents()

# ---------------------------------------------------------------------
# 4) Drop the now-stale hyperlink and shrink the conditional formatting
#    range down to the single still-populated id cell.
# ---------------------------------------------------------------------

$ws.Hyperlinks.Delete()

$fcs = $ws.Range("A2:A9").FormatConditions
$fcs.Item(1).ModifyAppliesToRange($ws.Range("A2"))

# ---------------------------------------------------------------------
# 5) Move the active selection.
# ---------------------------------------------------------------------

$null = $ws.Range("B10").Select()
